# Add the new "2022-Q3" sheet, positioned right before "2022-Q2",
# populate it with the Q3 fund-holding data, and insert the
# corresponding summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new worksheet right before the "2022-Q2" tab ---
$beforeSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q3"

# Copy the header-row / first-data-row styling from an existing
# quarterly sheet so the new tab matches the workbook's look.
$refSheet = $wb.Worksheets.Item(3)
$refSheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$refSheet.Range("A2").Copy($newSheet.Range("A2"))

# --- 2. Populate "2022-Q3" with its header + single data row ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").Value = "'516190"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "华夏中证文娱传媒ETF"
$newSheet.Range("D2").Value = "'0.13"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'96.01"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'2.75"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.0036"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 7

# --- 3. Shift the "总计" summary rows down and insert the Q3 row ---
$summary = $wb.Worksheets.Item(1)

# Extend the A-column "index" styling down onto the new last row
# before writing into it.
$summary.Range("A5").Copy($summary.Range("A6"))

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 25
$summary.Range("D6").Value = 5.16

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 4
$summary.Range("D5").Value = 1.27

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 13
$summary.Range("D4").Value = 4.32

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.03

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0

# --- 4. Restore original selection / active sheet ---
$summary.Activate()
$summary.Range("A1").Select()
